$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt11"
$ws.Range("C2").Value = "Fzd7"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 20.223983
$ws.Range("H2").Value = 60.671949
$ws.Range("I2").Value = 0.9624502889455165
$ws.Range("J2").Value = 0.9624502889455167
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.01111
$ws.Range("N2").Value = 3.03333
$ws.Range("O2").Value = 0.04063212692754557
$ws.Range("P2").Value = 0.04063212692754556
$ws.Range("Q2").Value = 20.44867145113
$ws.Range("R2").Value = 184.03804306017
$ws.Range("S2").Value = 0.03910640230188713
$ws.Range("T2").Value = 0.03910640230188713

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt11"
$ws.Range("C3").Value = "Fzd7"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 20.223983
$ws.Range("H3").Value = 60.671949
$ws.Range("I3").Value = 0.9624502889455165
$ws.Range("J3").Value = 0.9624502889455167
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.11799233333333
$ws.Range("N3").Value = 30.353977
$ws.Range("O3").Value = 0.4065982422683317
$ws.Range("P3").Value = 0.4065982422683317
$ws.Range("Q3").Value = 204.6261049434637
$ws.Range("R3").Value = 1841.634944491173
$ws.Range("S3").Value = 0.391330595755895
$ws.Range("T3").Value = 0.3913305957558951

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt11"
$ws.Range("C4").Value = "Fzd7"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 20.223983
$ws.Range("H4").Value = 60.671949
$ws.Range("I4").Value = 0.9624502889455165
$ws.Range("J4").Value = 0.9624502889455167
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 13.75539366666667
$ws.Range("N4").Value = 41.266181
$ws.Range("O4").Value = 0.5527696308041227
$ws.Range("P4").Value = 0.5527696308041226
$ws.Range("Q4").Value = 278.1888476729744
$ws.Range("R4").Value = 2503.699629056769
$ws.Range("S4").Value = 0.5320132908877344
$ws.Range("T4").Value = 0.5320132908877343

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Wnt11"
$ws.Range("C5").Value = "Fzd7"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7890326666666668
$ws.Range("H5").Value = 2.367098
$ws.Range("I5").Value = 0.03754971105448342
$ws.Range("J5").Value = 0.03754971105448343
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.01111
$ws.Range("N5").Value = 3.03333
$ws.Range("O5").Value = 0.04063212692754557
$ws.Range("P5").Value = 0.04063212692754556
$ws.Range("Q5").Value = 0.7977988195933334
$ws.Range("R5").Value = 7.18018937634
$ws.Range("S5").Value = 0.001525724625658431
$ws.Range("T5").Value = 0.001525724625658431

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Wnt11"
$ws.Range("C6").Value = "Fzd7"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7890326666666668
$ws.Range("H6").Value = 2.367098
$ws.Range("I6").Value = 0.03754971105448342
$ws.Range("J6").Value = 0.03754971105448343
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 10.11799233333333
$ws.Range("N6").Value = 30.353977
$ws.Range("O6").Value = 0.4065982422683317
$ws.Range("P6").Value = 0.4065982422683317
$ws.Range("Q6").Value = 7.983426472082891
$ws.Range("R6").Value = 71.85083824874602
$ws.Range("S6").Value = 0.01526764651243671
$ws.Range("T6").Value = 0.01526764651243671

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Wnt11"
$ws.Range("C7").Value = "Fzd7"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.7890326666666668
$ws.Range("H7").Value = 2.367098
$ws.Range("I7").Value = 0.03754971105448342
$ws.Range("J7").Value = 0.03754971105448343
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 13.75539366666667
$ws.Range("N7").Value = 41.266181
$ws.Range("O7").Value = 0.5527696308041227
$ws.Range("P7").Value = 0.5527696308041226
$ws.Range("Q7").Value = 10.85345494585978
$ws.Range("R7").Value = 97.68109451273803
$ws.Range("S7").Value = 0.02075633991638829
$ws.Range("T7").Value = 0.02075633991638829

# Remove now-unused rows 8-10 (data previously for MuSCs sending cluster rows are
# superseded by the shifted/recomputed rows above; delete the trailing rows so the
# table ends at row 7)
$ws.Range("A8:T10").Delete()

